$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D2:E51 are treated as text so numeric-looking strings (e.g. "0.561", "1.00")
# are stored as text rather than being coerced to numbers, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '61.464.34'
$ws.Range("E2").Value = '  -5.94%  '
$ws.Range("D3").Value = '2.981.77'
$ws.Range("E3").Value = '  -7.03%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '543.73'
$ws.Range("E5").Value = '  -5.53%  '
$ws.Range("D6").Value = '152.51'
$ws.Range("E6").Value = '  -9.21%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '0.561'
$ws.Range("E8").Value = '  -6.06%  '
$ws.Range("D9").Value = '2.984.49'
$ws.Range("E9").Value = '  -6.72%  '
$ws.Range("E10").Value = '  -6.71%  '
$ws.Range("D11").Value = '6.19'
$ws.Range("E11").Value = '  -8.32%  '
$ws.Range("E12").Value = '  -7.36%  '
$ws.Range("D13").Value = '3.501.01'
$ws.Range("E13").Value = '  -7.21%  '
$ws.Range("E14").Value = '  -3.66%  '
$ws.Range("D15").Value = '61.596.98'
$ws.Range("E15").Value = '  -5.68%  '
$ws.Range("D16").Value = '23.57'
$ws.Range("E16").Value = '  -8.34%  '
$ws.Range("D17").Value = '2.980.23'
$ws.Range("E17").Value = '  -6.88%  '
$ws.Range("D18").Value = '0.0000146'
$ws.Range("E18").Value = '  -7.35%  '
$ws.Range("D19").Value = '5.11'
$ws.Range("E19").Value = '  -4.50%  '
$ws.Range("D20").Value = '380.36'
$ws.Range("E20").Value = '  -8.04%  '
$ws.Range("D21").Value = '11.90'
$ws.Range("E21").Value = '  -7.90%  '
$ws.Range("D22").Value = '6.61'
$ws.Range("E22").Value = '  -8.20%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = '64.91'
$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D25").Value = '3.112.57'
$ws.Range("E25").Value = '  -7.00%  '
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").Value = '0.467'
$ws.Range("E26").Value = '  -5.07%  '
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = '0.186'
$ws.Range("E27").Value = '  -8.01%  '
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").Value = '0.995'
$ws.Range("E28").Value = '  -0.98%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0930'
$ws.Range("E29").Value = '  -12.31%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '8.29'
$ws.Range("E30").Value = '  -6.95%  '
$ws.Range("B31").Value = 'USDe'
$ws.Range("C31").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '1.72'
$ws.Range("E32").Value = '  -6.96%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '20.21'
$ws.Range("E33").Value = '  -6.31%  '
$ws.Range("B34").Value = 'Monero'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D34").Value = '158.20'
$ws.Range("E34").Value = '  +0.94%  '
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = '6.00'
$ws.Range("E35").Value = '  -6.67%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '4.60'
$ws.Range("E36").Value = '  -7.55%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").Value = '1.06'
$ws.Range("E37").Value = '  -7.26%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '1.27'
$ws.Range("E38").Value = '  -7.24%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '1.57'
$ws.Range("E39").Value = '  -8.92%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '2.430.36'
$ws.Range("E40").Value = '  -11.69%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '3.87'
$ws.Range("E41").Value = '  -7.02%  '
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = '37.00'
$ws.Range("E42").Value = '  -5.35%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '22.20'
$ws.Range("E43").Value = '  -8.93%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = '0.660'
$ws.Range("E44").Value = '  -7.62%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").Value = '0.0592'
$ws.Range("E45").Value = '  -6.61%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '0.0244'
$ws.Range("E47").Value = '  -7.37%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '4.93'
$ws.Range("E48").Value = '  -13.42%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '0.0953'
$ws.Range("E49").Value = '  -4.18%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '19.74'
$ws.Range("E50").Value = '  -8.96%  '
$ws.Range("B51").Value = 'WhiteBITCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D51").Value = '10.47'
$ws.Range("E51").Value = '  +0.22%  '

# Remove the temporary text formatting so cell styles match the original (no explicit
# style index), while the underlying stored values remain text.
$ws.Range("D2:E51").ClearFormats()

